$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that swap values between row 13 and row 14
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr13 = "$col" + "13"
    $addr14 = "$col" + "14"

    $val13 = $ws.Range($addr13).Value()
    $val14 = $ws.Range($addr14).Value()

    $ws.Range($addr13).Value = $val14
    $ws.Range($addr14).Value = $val13
}
